$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("上海电气", "永鼎股份", "山子高科"),
    @("永鼎股份", "上海电气", "上海电气"),
    @("山子高科", "江西铜业", "深科技"),
    @("深科技", "山子高科", "江西铜业"),
    @("江西铜业", "深科技", "杉杉股份"),
    @("合锻智能", "浪潮信息", "永鼎股份"),
    @("浪潮信息", "东方财富", "合锻智能"),
    @("杉杉股份", "合锻智能", "三花智控"),
    @("大众公用", "大众公用", "万向钱潮"),
    @("三花智控", "贵州茅台", "浪潮信息"),
    @("华友钴业", "华友钴业", "天赐材料"),
    @("长电科技", "杉杉股份", "华友钴业"),
    @("天赐材料", "江波龙", "中电鑫龙"),
    @("特变电工", "兆易创新", "海南华铁"),
    @("兆易创新", "长电科技", "张江高科"),
    @("东方财富", "特变电工", "大众公用"),
    @("中芯国际", "三花智控", "特变电工"),
    @("阳光电源", "赣锋锂业", "兆易创新"),
    @("江波龙", "洛阳钼业", "长电科技"),
    @("寒武纪-U", "云南铜业", "洛阳钼业")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
